$wb = $excel.ActiveWorkbook

# The "user" table worksheet
$ws = $wb.Worksheets.Item("user")

# Insert a new row at position 5 (pushing existing rows 5+ down by one),
# matching the gap pattern already present in the sheet.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with the new "level" column definition.
$ws.Range("A5").Value = "level"
$ws.Range("B5").Value = "1：管理员 2：普通用户"

# Make the "user" sheet the active/selected sheet and move the selection
# to the first empty cell below the table data (mirrors the gap pattern
# where row 11, i.e. one row below the last populated row 10, stays empty).
$ws.Activate()
$ws.Range("A13").Select()
